$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max") - this shifts D->C, E->D
$ws.Range("C1:C1").EntireColumn.Delete()

# Delete row 3 (even_MAG-GUT14745.fa)
$ws.Range("A3:A3").EntireRow.Delete()

# Update B2 value
$ws.Range("B2").Value = 139.0436161570187
